$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly report gained one new week of data (2022-03-24, serial 44644).
# Two new detail rows are inserted right before the old row 210, pushing the
# rest of the table (old rows 210-255) down by two rows without altering any
# of their values.
$ws.Rows("210:211").Insert()

# New row 210: "Primera" quality entry for the new week, mirroring the
# existing pattern already present for this market/date group.
$ws.Cells.Item(210,1).Value = 9
$ws.Cells.Item(210,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(210,3).Value = "Metropolitana"
$ws.Cells.Item(210,4).Value = 44644
$ws.Cells.Item(210,5).Value = 13
$ws.Cells.Item(210,6).Value = 100112017
$ws.Cells.Item(210,7).Value = "Apio"
$ws.Cells.Item(210,8).Value = "Americana (o)"
$ws.Cells.Item(210,9).Value = "Primera"
$ws.Cells.Item(210,10).Value = 79
$ws.Cells.Item(210,11).Value = 9000
$ws.Cells.Item(210,12).Value = 10000
$ws.Cells.Item(210,13).Value = 9494
$ws.Cells.Item(210,14).Value = "`$/docena de matas"
$ws.Cells.Item(210,15).Value = "Región de Coquimbo"
$ws.Cells.Item(210,16).Value = 1582
$ws.Cells.Item(210,17).Value = 6
$ws.Cells.Item(210,18).Value = "Hortaliza"

# New row 211: "Segunda" quality entry for the same new week.
$ws.Cells.Item(211,1).Value = 9
$ws.Cells.Item(211,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(211,3).Value = "Metropolitana"
$ws.Cells.Item(211,4).Value = 44644
$ws.Cells.Item(211,5).Value = 13
$ws.Cells.Item(211,6).Value = 100112017
$ws.Cells.Item(211,7).Value = "Apio"
$ws.Cells.Item(211,8).Value = "Americana (o)"
$ws.Cells.Item(211,9).Value = "Segunda"
$ws.Cells.Item(211,10).Value = 43
$ws.Cells.Item(211,11).Value = 8000
$ws.Cells.Item(211,12).Value = 8000
$ws.Cells.Item(211,13).Value = 8000
$ws.Cells.Item(211,14).Value = "`$/docena de matas"
$ws.Cells.Item(211,15).Value = "Región de Coquimbo"
$ws.Cells.Item(211,16).Value = 1333
$ws.Cells.Item(211,17).Value = 6
$ws.Cells.Item(211,18).Value = "Hortaliza"
